# improve sensitivity analysis plots
#
# - add a new "immune_cells" worksheet (color lookup for a sensitivity-
#   analysis plot restricted to the immune-cell subset) after "cell_types"
# - re-colour the "methods" sheet using the ColorBrewer "Dark2" palette
#   (and drop the now-unused "backup2" row)
# - leave "cell_types" values as-is

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. add the new "immune_cells" sheet after "cell_types"
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$immune = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$immune.Name = "immune_cells"

# fill column A (value labels) first ...
$immune.Range("A1").Value = "value"
$immune.Range("A2").Value = "B cell"
$immune.Range("A3").Value = "Dendritic cell"
$immune.Range("A4").Value = "Macrophage/Monocyte"
$immune.Range("A5").Value = "NK cell"
$immune.Range("A6").Value = "T cell CD4+"
$immune.Range("A7").Value = "T cell CD8+"
$immune.Range("A8").Value = "T cell CD4+ (non-regulatory)"
$immune.Range("A9").Value = "T cell regulatory (Tregs)"

# ... then column B (Dark2 palette colors)
$immune.Range("B1").Value = "color"
$immune.Range("B2").Value = "#1b9e77"
$immune.Range("B3").Value = "#d95f02"
$immune.Range("B4").Value = "#7570b3"
$immune.Range("B5").Value = "#e7298a"
$immune.Range("B6").Value = "#66a61e"
$immune.Range("B7").Value = "#e6ab02"
$immune.Range("B8").Value = "#a6761d"
$immune.Range("B9").Value = "#666666"

$immune.Columns.Item(1).ColumnWidth = 26.140625
$immune.Range("B9").Select() | Out-Null

# ---------------------------------------------------------------------
# 2. "methods" sheet: recolor with the same Dark2 palette, drop the
#    now-unused backup2 row
# ---------------------------------------------------------------------
$methods = $wb.Worksheets.Item("methods")

$methods.Range("B2").Value = "#1b9e77"
$methods.Range("B3").Value = "#d95f02"
$methods.Range("B4").Value = "#7570b3"
$methods.Range("B5").Value = "#e7298a"
$methods.Range("B6").Value = "#66a61e"
$methods.Range("B7").Value = "#e6ab02"
$methods.Range("B8").Value = "#a6761d"
$methods.Range("B9").Value = "#666666"

# row 10 held "backup2" -- no longer used, remove it entirely
$methods.Rows.Item(10).Delete()

# ---------------------------------------------------------------------
# 3. "cell_types" sheet: content is unchanged, just move the selection
#    and make it inactive
# ---------------------------------------------------------------------
$cellTypes = $wb.Worksheets.Item("cell_types")
$cellTypes.Range("A15").Select() | Out-Null

# ---------------------------------------------------------------------
# 4. make "methods" the active sheet/tab again, with A10:B10 selected
# ---------------------------------------------------------------------
$methods.Range("A10:B10").Select() | Out-Null
